$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'69.133.01"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -1.98%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'3.519.56"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -2.94%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.27%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'583.17"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -3.42%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'192.47"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -2.04%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.605"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  -3.59%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D9').Value = "'0.203"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -2.36%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.617"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -4.65%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'52.04"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -2.45%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'0.0000285"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -5.76%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'9.15"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -4.38%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'4.083.45"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -2.85%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'643.08"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +6.98%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'69.256.66"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -1.96%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'3.539.22"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -2.85%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'12.45"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -4.08%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'  -1.71%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'18.24"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -4.23%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'0.952"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -4.59%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'17.93"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -2.05%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = "'  +4.76%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'101.57"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -0.58%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'4.34"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -5.94%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'2.89"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -3.59%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'10.06"
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Value = "'9.43"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -2.61%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'32.76"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -3.20%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('B30').Value = "'NEARProtocol"
$ws.Range('B30').Style = 'Normal'
$ws.Range('C30').Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range('C30').Style = 'Normal'
$ws.Range('D30').Value = "'6.71"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -8.07%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('B31').Value = "'dogwifhat"
$ws.Range('B31').Style = 'Normal'
$ws.Range('C31').Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range('C31').Style = 'Normal'
$ws.Range('D31').Value = "'4.08"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -12.35%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'11.63"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  -5.34%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'0.109"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -6.99%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'61.29"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -3.53%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'3.712.00"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -5.44%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'0.998"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -0.33%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'0.0₃0794"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -10.67%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'504.08"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -5.88%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'3.57"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +0.71%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'2.93"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -4.36%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'0.367"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -5.77%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('B42').Value = "'InjectiveProtocol"
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').Value = "'34.52"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -6.50%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('B43').Value = "'Kaspa"
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Value = "'0.133"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -0.70%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'0.0441"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -4.54%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'3.39"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -0.98%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'2.83"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -1.60%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = "'  -4.17%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'  -0.05%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'8.17"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -4.80%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'2.70"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +58.68%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'0.153"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +1.34%  "
$ws.Range('E51').Style = 'Normal'
